# Adds a new "Images" worksheet (Hashtags / Filename table) as the last tab
# of the workbook, formats its header row bold, sizes its columns, updates
# the Issues sheet's remembered selection to B6, and finally leaves the new
# Images sheet as the active tab - matching the authored commit.

$wb = $excel.ActiveWorkbook

# --- add the new "Images" sheet as the LAST tab ----------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$images = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$images.Name = "Images"

# --- header row (bold) ------------------------------------------------------
$images.Range("A1").Value = "Hashtags"
$images.Range("B1").Value = "Filename"
$images.Range("A1:B1").Font.Bold = $true

# --- data row ----------------------------------------------------------------
$images.Range("A2").Value = "#Klima #CO2 #Umweltschutz"
$images.Range("B2").Value = "windmills-5643293_1280.jpg"

# --- column widths -------------------------------------------------------------
$images.Columns.Item(1).ColumnWidth = 31.140625
$images.Columns.Item(2).ColumnWidth = 37.140625

# --- page setup (paper size / orientation / margins) ---------------------------
$images.PageSetup.PaperSize = 9
$images.PageSetup.Orientation = 1
$images.PageSetup.LeftMargin = $excel.InchesToPoints(0.7)
$images.PageSetup.RightMargin = $excel.InchesToPoints(0.7)
$images.PageSetup.TopMargin = $excel.InchesToPoints(2 / 2.54)
$images.PageSetup.BottomMargin = $excel.InchesToPoints(2 / 2.54)
$images.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$images.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)

# --- restore/update the Issues sheet selection (B6) -----------------------------
$issues = $wb.Worksheets.Item("Issues")
$issues.Activate()
[void]$issues.Range("B6").Select()

# --- make the new Images sheet the active tab -----------------------------------
$images.Activate()
